# Move regression estimates to country specific folder.
# The workbook's single sheet, previously named "EL" (doubling as the
# country-code folder name), is renamed to "Parameters" now that the
# file itself lives inside the country-specific "EL" folder.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Parameters"
